# Apply the "Add files via upload" banner edit to the single slide:
#   1. Remove the first banner group ("Group 21") entirely - it duplicated
#      the same picture/overlay/text as the other banner.
#   2. In the remaining banner group ("Group 1"), update the date textbox's
#      second line to read "Lyon, 15-17 February 2017" instead of
#      "15-17 February 2017".
#   3. Remove the stray standalone picture ("Picture 11") and textbox
#      ("TextBox 13") that duplicated the banner content outside any group.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# 1. Delete "Group 21" (first shape in the slide).
$group21 = $s.Shapes("Group 21")
$group21.Delete()

# 2. Update the text of "TextBox 9" inside "Group 1".
$group1 = $s.Shapes("Group 1")
$textBox9 = $group1.GroupItems("TextBox 9")
$dateRun = $textBox9.TextFrame.TextRange.Paragraphs(2, 1).Runs(1, 1)
$dateRun.Text = "Lyon, 15-17 February 2017"

# 3. Delete the standalone "Picture 11" and "TextBox 13" shapes.
$s.Shapes("Picture 11").Delete()
$s.Shapes("TextBox 13").Delete()
